$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Concept" column (B) with a "Text" column.
$ws.Range("B1").Value = "Text"
$ws.Range("B2").Value = "hypotheses"
$ws.Range("B3").Value = "hypotheses"
$ws.Range("B4").Value = "hypotheses"
$ws.Range("B5").Value = "hypotheses"

# Update the active selection to match the saved view state.
$ws.Range("B6").Select()
